# M10: Add hotmail address to excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: 번호=4, Hosting업체=MS, 이메일 주소=iskim0706@hotmail.com, PS=dsfsf
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "MS"
$ws.Range("D6").Value = "iskim0706@hotmail.com"
$ws.Range("F6").Value = "dsfsf"

# Add mailto hyperlink on the new email cell, matching the style used by the
# existing hyperlink cell (D3)
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:iskim0706@hotmail.com")
$ws.Range("D6").Style = $ws.Range("D3").Style

# Update the active selection to match the saved view state
$ws.Range("G12").Select()
